# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns for both language sheets (the handback data that was missing before),
# updates the Status text to reflect the handback, and widens a couple of columns
# that now need to show longer file names.

$wb = $excel.ActiveWorkbook

$mdUrl0351 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb9251b23a89679a096d069d7da044ee1344f02e/e2e/0351e35a-635a-48de-af8c-b4ade49b12f0.md"
$mdUrlCbff = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb9251b23a89679a096d069d7da044ee1344f02e/e2e/cbff6d1b-797e-477b-98a1-18ce611332a3.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Overview sheet: update the per-language status cells (columns E = zh-cn, F = de-de)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

# Rebuild the hyperlinks collection so the new "Latest Target File" links land
# in the same relationship-id order Excel would generate (A2, I2, A3, I3).
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl0351, "", "", "0351e35a-635a-48de-af8c-b4ade49b12f0.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl0351, "", "", "0351e35a-635a-48de-af8c-b4ade49b12f0.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrlCbff, "", "", "cbff6d1b-797e-477b-98a1-18ce611332a3.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrlCbff, "", "", "cbff6d1b-797e-477b-98a1-18ce611332a3.md")

$wsZh.Range("J2").Value = "0351e35a-635a-48de-af8c-b4ade49b12f0.06d3aa00ef44f790689b00f44dc1bc837782aee9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-11-03 20:25:33"
$wsZh.Range("J3").Value = "cbff6d1b-797e-477b-98a1-18ce611332a3.583f98a3a62b8660e75a4ed7d9a63c244083fc28.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-11-03 20:25:33"

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl0351, "", "", "0351e35a-635a-48de-af8c-b4ade49b12f0.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl0351, "", "", "0351e35a-635a-48de-af8c-b4ade49b12f0.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrlCbff, "", "", "cbff6d1b-797e-477b-98a1-18ce611332a3.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrlCbff, "", "", "cbff6d1b-797e-477b-98a1-18ce611332a3.md")

$wsDe.Range("J2").Value = "0351e35a-635a-48de-af8c-b4ade49b12f0.06d3aa00ef44f790689b00f44dc1bc837782aee9.de-de.xlf"
$wsDe.Range("K2").Value = "2016-11-03 20:25:51"
$wsDe.Range("J3").Value = "cbff6d1b-797e-477b-98a1-18ce611332a3.583f98a3a62b8660e75a4ed7d9a63c244083fc28.de-de.xlf"
$wsDe.Range("K3").Value = "2016-11-03 20:25:51"

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
